# Auto-generated Excel COM-interop script
# Refresh the cryptos table (rows 2-51, columns B:E):
#  - updated Price / Volume(1h) figures for every existing coin
#  - a new coin ("OKB") inserted at row 9, pushing every coin that was
#    below it down by one row
#  - the previously-last coin ("Aave") drops off the bottom of the
#    fixed 50-row table as a result of that shift

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a price string to column D (Price) while keeping it as
# literal text (Excel would otherwise parse "46.90" / "1.004" as a number
# and normalize away the formatting the source data relies on). The
# NumberFormat flip forces a text entry; resetting the style back to
# "Normal" afterwards avoids leaving a stray text-format style behind.
function Set-PriceText($row, $text) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "29.146.77", "  +1.68%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.908.44", "  +1.95%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.004", "  -0.23%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "327.50", "  +0.20%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.003", "  -0.27%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4647", "  -0.50%  ")
    ,@(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3930", "  +1.08%  ")
    ,@(9, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "46.90", "  +0.96%  ")
    ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07971", "  +1.18%  ")
    ,@(11, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.005", "  +3.15%  ")
    ,@(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "22.38", "  +1.70%  ")
    ,@(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.943.35", "  +3.52%  ")
    ,@(14, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.155", "  +2.22%  ")
    ,@(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.812", "  +1.77%  ")
    ,@(16, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06995", "  +0.13%  ")
    ,@(17, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "88.62", "  +0.49%  ")
    ,@(18, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.004", "  -0.35%  ")
    ,@(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001011", "  +0.58%  ")
    ,@(20, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "17.27", "  +2.64%  ")
    ,@(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.002", "  -0.35%  ")
    ,@(22, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "29.151.35", "  +1.70%  ")
    ,@(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.390", "  +1.65%  ")
    ,@(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.07", "  +0.45%  ")
    ,@(25, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.135.23", "  +2.92%  ")
    ,@(26, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.054", "  -3.01%  ")
    ,@(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "155.88", "  +1.94%  ")
    ,@(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "19.60", "  +1.87%  ")
    ,@(29, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.864", "  +2.06%  ")
    ,@(30, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.014", "  +1.26%  ")
    ,@(31, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "119.71", "  +0.29%  ")
    ,@(32, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09398", "  +0.26%  ")
    ,@(33, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.9291", "  +0.94%  ")
    ,@(34, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.358", "  +1.57%  ")
    ,@(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.344", "  +0.26%  ")
    ,@(36, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.279", "  -1.95%  ")
    ,@(37, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05859", "  +0.74%  ")
    ,@(38, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "8.036", "  +3.49%  ")
    ,@(39, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.160", "  +1.22%  ")
    ,@(40, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02107", "  +0.23%  ")
    ,@(41, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.5766", "  +2.45%  ")
    ,@(42, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1812", "  +1.36%  ")
    ,@(43, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "10.00", "  +2.42%  ")
    ,@(44, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "12.06", "  +2.98%  ")
    ,@(45, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.241", "  +8.55%  ")
    ,@(46, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.5442", "  +2.37%  ")
    ,@(47, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.07113", "  -1.49%  ")
    ,@(48, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.887", "  +3.36%  ")
    ,@(49, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.118", "  -3.42%  ")
    ,@(50, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.568", "  +6.41%  ")
    ,@(51, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "112.59", "  -0.53%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    Set-PriceText $r $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

